# Remove the now-obsolete "© 2020 ..." footer block (the blank paragraph,
# the page-break paragraph, and the copyright paragraph that preceded it)
# that used to sit right after the "LOT2039: ... (Requisito)" paragraph.
# A duplicate blank + page-break paragraph pair already follows, so the
# trailing structure of the document stays intact.

$d = $word.ActiveDocument

$target = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$copyrightPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $copyrightPara = $p
        break
    }
}

if ($copyrightPara -ne $null) {
    $pageBreakPara = $copyrightPara.Previous()
    $blankPara = $pageBreakPara.Previous()

    $start = $blankPara.Range.Start
    $end = $copyrightPara.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
